# Sync attendance_reports: reorder "Recorded By" names so that
# "System"/"system" entries are moved to the front of the comma-separated
# list (i.e. the whole list is reversed) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column G ("Recorded By")
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ",\s*"

    if ($parts.Count -gt 1) {
        # Determine whether one of the names is the "System" account and
        # whether it is already the first entry.
        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p -eq "System" -or $p -eq "system") {
                $hasSystem = $true
            }
        }
        $firstIsSystem = ($parts[0] -eq "System" -or $parts[0] -eq "system")

        if ($hasSystem -and -not $firstIsSystem) {
            $reversed = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
